$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "Author") {
        $p.Range.Delete()
        break
    }
}
